# Applies the "added ifoCAST full series evaluation" update:
# refreshes the naive QoQ error-series values in B2:K20 with the latest
# ifoCAST-matched evaluation numbers, and trims the trailing (oldest-vintage)
# cells off rows 11-20 to match the new, shorter forecast horizons.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$rowsData = @(
    @(-2.808131625772977, -0.77139249106923, -0.1687823066887937, -0.4024803045374513, 0.3742746160819611, -0.3294423225852309, 0.1540433631259821, 0.6252233516998879, 0.9388562853074672, 0.0600786799400978),
    @(-0.3631600374474599, 0.08284073719760038, -0.4843863583605866, 0.4493666654632467, -0.2235351620715071, 0.2124600092846375, 0.6916620355299794, 1.013538243918304, 0.1303919718410766, -0.5999993966742521),
    @(-0.02265847465121479, -0.6104916423259444, 0.6176147340964784, -0.1628710694676248, 0.2235928486592704, 0.75452478685691, 1.070889771712047, 0.1748971488454556, -0.5489255256478335, 1.115476138671014),
    @(-1.214932262173683, 0.553411356741785, 0.128676238589592, 0.1858859725733368, 0.7459407667834489, 1.154630849309415, 0.2097242937735006, -0.5223611425955365, 1.160974979417613, 0.5338127914751168),
    @(0.4657705279028737, 0.07625703655992344, 0.2445323695583424, 0.7508699121205502, 1.14353756626052, 0.2224909929985871, -0.5136315835481774, 1.163476923566464, 0.540088728946798, 0.7747069251866952),
    @(0.09868501619293202, 0.2471359091207911, 0.728802226466535, 1.137555585956165, 0.217147742816176, -0.5240001169185436, 1.154999601442841, 0.5323164206399872, 0.7660829516537266, 0.6652412546050546),
    @(0.2685604805780112, 0.8489800752657474, 1.048361241542382, 0.1874207389675982, -0.5153739968886335, 1.129059428017406, 0.5095378576117597, 0.7524649165051087, 0.6467204574198788, -0.1454929044188731),
    @(0.806659442945358, 1.019156767686649, 0.184400139162677, -0.534235280500114, 1.108728296952201, 0.4950077021775119, 0.7358238710128409, 0.6291366119911286, -0.1619667879489159, 0.449372724506711),
    @(1.359354508304559, 0.2597712009466141, -0.7253285668131131, 1.133422209538934, 0.5043871057338079, 0.6810130724014498, 0.6102425211024588, -0.1768307887639616, 0.4206610138773402, 0.2638965897873631),
    @(0.7104660729368646, -0.6782178646007546, 0.8931389787481262, 0.5403321867996707, 0.6789623149618569, 0.5336742816664286, -0.204230183139245, 0.3943719636796149, 0.2204992990740305, $null),
    @(-0.3658922776772162, 1.026427718483651, 0.3612537649834791, 0.707156553732424, 0.5696995928588221, -0.2417717533434586, 0.3888296245922537, 0.2246746280127792, $null, $null),
    @(1.262018209591492, 0.4459325058577887, 0.5615561870987069, 0.5796533357180647, -0.2264065836439137, 0.3499812896348306, 0.2088288189855932, $null, $null, $null),
    @(0.7569566923391715, 0.6786449615099022, 0.4171891942684979, -0.1974476331787121, 0.3857241620897341, 0.1775011726019661, $null, $null, $null, $null),
    @(0.9254701389140165, 0.4348450618063874, -0.2920698722897066, 0.4179415503382142, 0.1898892984296834, $null, $null, $null, $null, $null),
    @(0.6745214212225993, -0.2070456288204931, 0.3014569719802002, 0.2049945700815359, $null, $null, $null, $null, $null, $null),
    @(-0.04218555178640582, 0.3149942442281164, 0.1420216510915729, $null, $null, $null, $null, $null, $null, $null),
    @(0.5688432860935244, 0.2413397012736094, $null, $null, $null, $null, $null, $null, $null, $null),
    @(0.2853993925130583, $null, $null, $null, $null, $null, $null, $null, $null, $null),
    @($null, $null, $null, $null, $null, $null, $null, $null, $null, $null)
)

$startRow = 2
$startCol = 2  # column B

for ($i = 0; $i -lt $rowsData.Count; $i++) {
    $row = $rowsData[$i]
    for ($j = 0; $j -lt $row.Count; $j++) {
        $cell = $ws.Cells.Item($startRow + $i, $startCol + $j)
        $val = $row[$j]
        if ($null -eq $val) {
            $cell.ClearContents()
        } else {
            $cell.Value = $val
        }
    }
}
